# Auto-generated edit script applying the scheduled-runner update to Titan_Profits workbook.
# For each touched cell: write the new numeric value, or clear the cell when the diff
# removes it entirely (no replacement <v>).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 601.5625
$ws.Range("I19").Value = 481.33334
$ws.Range("K19").Value = 481.33334
$ws.Range("M19").Value = -306.33334
$ws.Range("H100").Value = 16669027
$ws.Range("I100").Value = 20835034
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 20835034
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -20834493
$ws.Range("N100").Value = -6082
$ws.Range("H133").Value = 46046
$ws.Range("J133").Value = 46046
$ws.Range("L133").Value = 46046
$ws.Range("N133").Value = -56166
$ws.Range("H138").Value = 13495276
$ws.Range("I138").Value = 2301410.2
$ws.Range("J138").Value = 38466210
$ws.Range("K138").Value = 6904230.600000001
$ws.Range("L138").Value = 115398630
$ws.Range("M138").Value = -6899090.600000001
$ws.Range("N138").Value = -115408910
$ws.Range("H141").Value = 4414.423
$ws.Range("I141").Value = 2738.9473
$ws.Range("J141").Value = 8962.143
$ws.Range("K141").Value = 8216.841899999999
$ws.Range("L141").Value = 26886.429
$ws.Range("M141").Value = -3036.841899999999
$ws.Range("N141").Value = -37246.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 16250
$ws.Range("J9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("N9").Value = -5340
$ws.Range("H20").Value = 16250
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5540
$ws.Range("H23").Value = 10000000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H37").Value = 4300
$ws.Range("J37").Value = 4300
$ws.Range("L37").Value = 4300
$ws.Range("N37").Value = -4846
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 7400
$ws.Range("J55").Value = 7400
$ws.Range("L55").Value = 7400
$ws.Range("N55").Value = -8030
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240
$ws.Range("H80").Value = 75815.5
$ws.Range("I80").Value = 99000
$ws.Range("J80").Value = 52631
$ws.Range("K80").Value = 99000
$ws.Range("L80").Value = 52631
$ws.Range("M80").Value = -98002
$ws.Range("N80").Value = -54627
$ws.Range("H83").Value = 75815.5
$ws.Range("I83").Value = 99000
$ws.Range("J83").Value = 52631
$ws.Range("K83").Value = 297000
$ws.Range("L83").Value = 157893
$ws.Range("M83").Value = -292008
$ws.Range("N83").Value = -167877
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H102").Value = 1765.909
$ws.Range("I102").Value = 1579.3334
$ws.Range("J102").Value = 2605.5
$ws.Range("K102").Value = 1579.3334
$ws.Range("L102").Value = 2605.5
$ws.Range("M102").Value = 42.66660000000002
$ws.Range("N102").Value = -5849.5
$ws.Range("H132").Value = 3133.3333
$ws.Range("I132").Value = 2723.111
$ws.Range("J132").Value = 4364
$ws.Range("K132").Value = 8169.333
$ws.Range("L132").Value = 13092
$ws.Range("M132").Value = -5639.333
$ws.Range("N132").Value = -18152
$ws.Range("H133").Value = 46904.4
$ws.Range("J133").Value = 46904.4
$ws.Range("L133").Value = 46904.4
$ws.Range("N133").Value = -51964.4
$ws.Range("H139").Value = 56238.332
$ws.Range("J139").Value = 56238.332
$ws.Range("L139").Value = 56238.332
$ws.Range("N139").Value = -66518.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 670.8387
$ws.Range("I94").Value = 522.0833
$ws.Range("K94").Value = 522.0833
$ws.Range("M94").Value = -71.08330000000001
$ws.Range("H105").Value = 3292.2327
$ws.Range("I105").Value = 3072.742
$ws.Range("J105").Value = 3859.25
$ws.Range("K105").Value = 3072.742
$ws.Range("L105").Value = 3859.25
$ws.Range("M105").Value = -1325.742
$ws.Range("N105").Value = -7353.25
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H126").Value = 56333.332
$ws.Range("I126").Value = 56000
$ws.Range("J126").Value = 56500
$ws.Range("K126").Value = 56000
$ws.Range("L126").Value = 56500
$ws.Range("M126").Value = -51060
$ws.Range("N126").Value = -66380
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H139").Value = 94793.336
$ws.Range("J139").Value = 94793.336
$ws.Range("L139").Value = 94793.336
$ws.Range("N139").Value = -105073.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3501.3062
$ws.Range("I31").Value = 1519.6
$ws.Range("K31").Value = 1519.6
$ws.Range("M31").Value = -1224.6
$ws.Range("H34").Value = 3501.3062
$ws.Range("I34").Value = 1519.6
$ws.Range("K34").Value = 1519.6
$ws.Range("M34").Value = -1317.6
$ws.Range("H41").Value = 2000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H99").Value = 2557.4736
$ws.Range("I99").Value = 1924.3334
$ws.Range("J99").Value = 3642.8572
$ws.Range("K99").Value = 1924.3334
$ws.Range("L99").Value = 3642.8572
$ws.Range("M99").Value = -426.3334
$ws.Range("N99").Value = -6638.8572
$ws.Range("H126").Value = 2557.4736
$ws.Range("I126").Value = 1924.3334
$ws.Range("J126").Value = 3642.8572
$ws.Range("K126").Value = 5773.0002
$ws.Range("L126").Value = 10928.5716
$ws.Range("M126").Value = -3303.0002
$ws.Range("N126").Value = -15868.5716
$ws.Range("H132").Value = 1878.4839
$ws.Range("I132").Value = 1363.3043
$ws.Range("J132").Value = 3359.625
$ws.Range("K132").Value = 4089.9129
$ws.Range("L132").Value = 10078.875
$ws.Range("M132").Value = -1559.9129
$ws.Range("N132").Value = -15138.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 267.25
$ws.Range("I47").Value = 267.25
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 801.75
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -370.75
$ws.Range("N47").ClearContents()
$ws.Range("H131").Value = 11496320
$ws.Range("I131").Value = 750
$ws.Range("J131").Value = 12347844
$ws.Range("K131").Value = 2250
$ws.Range("L131").Value = 37043532
$ws.Range("M131").Value = 2790
$ws.Range("N131").Value = -37053612

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9000
$ws.Range("J92").Value = 9000
$ws.Range("L92").Value = 9000
$ws.Range("N92").Value = -12744
$ws.Range("H122").Value = 3203.6667
$ws.Range("I122").Value = 1207
$ws.Range("J122").Value = 3603
$ws.Range("K122").Value = 3621
$ws.Range("L122").Value = 10809
$ws.Range("M122").Value = -1171
$ws.Range("N122").Value = -15709
$ws.Range("H132").Value = 4019.1924
$ws.Range("I132").Value = 3357.7368
$ws.Range("J132").Value = 5814.5713
$ws.Range("K132").Value = 10073.2104
$ws.Range("L132").Value = 17443.7139
$ws.Range("M132").Value = -7543.2104
$ws.Range("N132").Value = -22503.7139
$ws.Range("H137").Value = 60390
$ws.Range("J137").Value = 60390
$ws.Range("L137").Value = 60390
$ws.Range("N137").Value = -70590
$ws.Range("H138").Value = 74428.625
$ws.Range("J138").Value = 74428.625
$ws.Range("L138").Value = 74428.625
$ws.Range("N138").Value = -84708.625
$ws.Range("H139").Value = 47000
$ws.Range("J139").Value = 47000
$ws.Range("L139").Value = 47000
$ws.Range("N139").Value = -57280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12501687
$ws.Range("I16").Value = 16668600
$ws.Range("J16").Value = 950
$ws.Range("K16").Value = 16668600
$ws.Range("L16").Value = 950
$ws.Range("M16").Value = -16668430
$ws.Range("N16").Value = -1290
$ws.Range("H46").Value = 1302.7142
$ws.Range("I46").Value = 806.6667
$ws.Range("J46").Value = 1674.75
$ws.Range("K46").Value = 806.6667
$ws.Range("L46").Value = 1674.75
$ws.Range("M46").Value = -618.6667
$ws.Range("N46").Value = -2050.75
$ws.Range("H93").Value = 3222.875
$ws.Range("I93").Value = 2547.1667
$ws.Range("J93").Value = 5250
$ws.Range("K93").Value = 2547.1667
$ws.Range("L93").Value = 5250
$ws.Range("M93").Value = -1299.1667
$ws.Range("N93").Value = -7746
$ws.Range("H94").Value = 49000
$ws.Range("J94").Value = 49000
$ws.Range("L94").Value = 49000
$ws.Range("N94").Value = -50352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 33316.668
$ws.Range("J63").Value = 33316.668
$ws.Range("L63").Value = 33316.668
$ws.Range("N63").Value = -34564.668
$ws.Range("H66").Value = 33316.668
$ws.Range("J66").Value = 33316.668
$ws.Range("L66").Value = 99950.00399999999
$ws.Range("N66").Value = -106190.004
$ws.Range("H86").Value = 50000
$ws.Range("J86").Value = 50000
$ws.Range("L86").Value = 50000
$ws.Range("N86").Value = -52246
$ws.Range("H89").Value = 50000
$ws.Range("J89").Value = 50000
$ws.Range("L89").Value = 250000
$ws.Range("N89").Value = -261232
$ws.Range("H110").Value = 39640
$ws.Range("J110").Value = 39640
$ws.Range("L110").Value = 39640
$ws.Range("N110").Value = -47820
$ws.Range("H132").Value = 1422.3934
$ws.Range("I132").Value = 1106.5
$ws.Range("J132").Value = 3247.5557
$ws.Range("K132").Value = 3319.5
$ws.Range("L132").Value = 9742.667099999999
$ws.Range("M132").Value = -789.5
$ws.Range("N132").Value = -14802.6671

